$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range('D2') '62.869.80'
Set-TextValue $ws.Range('E2') '  +2.22%  '
Set-TextValue $ws.Range('D3') '3.041.02'
Set-TextValue $ws.Range('E3') '  +1.50%  '
Set-TextValue $ws.Range('E4') '  +0.05%  '
Set-TextValue $ws.Range('D5') '595.93'
Set-TextValue $ws.Range('E5') '  +1.45%  '
Set-TextValue $ws.Range('D6') '153.73'
Set-TextValue $ws.Range('E6') '  +6.64%  '
Set-TextValue $ws.Range('E7') '  -0.01%  '
Set-TextValue $ws.Range('D8') '3.034.38'
Set-TextValue $ws.Range('E8') '  +1.31%  '
Set-TextValue $ws.Range('E9') '  -0.34%  '
Set-TextValue $ws.Range('D10') '6.34'
Set-TextValue $ws.Range('E10') '  +7.37%  '
Set-TextValue $ws.Range('E11') '  +3.07%  '
Set-TextValue $ws.Range('D12') '0.465'
Set-TextValue $ws.Range('E12') '  +0.52%  '
Set-TextValue $ws.Range('D13') '0.0000235'
Set-TextValue $ws.Range('E13') '  +3.22%  '
Set-TextValue $ws.Range('D14') '35.35'
Set-TextValue $ws.Range('E14') '  +2.86%  '
Set-TextValue $ws.Range('E15') '  +2.08%  '
Set-TextValue $ws.Range('D16') '3.544.13'
Set-TextValue $ws.Range('E16') '  +1.63%  '
Set-TextValue $ws.Range('B17') 'Polkadot'
Set-TextValue $ws.Range('C17') 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue $ws.Range('D17') '7.08'
Set-TextValue $ws.Range('E17') '  +0.69%  '
Set-TextValue $ws.Range('B18') 'WrappedBTC'
Set-TextValue $ws.Range('C18') 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue $ws.Range('D18') '62.873.47'
Set-TextValue $ws.Range('E18') '  +2.34%  '
Set-TextValue $ws.Range('D19') '3.037.46'
Set-TextValue $ws.Range('E19') '  +1.51%  '
Set-TextValue $ws.Range('D20') '452.82'
Set-TextValue $ws.Range('E20') '  -0.10%  '
Set-TextValue $ws.Range('D21') '14.30'
Set-TextValue $ws.Range('E21') '  +1.68%  '
Set-TextValue $ws.Range('D22') '0.696'
Set-TextValue $ws.Range('E22') '  +1.13%  '
Set-TextValue $ws.Range('D23') '7.52'
Set-TextValue $ws.Range('E23') '  +2.22%  '
Set-TextValue $ws.Range('D24') '83.36'
Set-TextValue $ws.Range('E24') '  +1.74%  '
Set-TextValue $ws.Range('D25') '2.31'
Set-TextValue $ws.Range('E25') '  +5.69%  '
Set-TextValue $ws.Range('D26') '11.00'
Set-TextValue $ws.Range('E26') '  +6.45%  '
Set-TextValue $ws.Range('D27') '12.26'
Set-TextValue $ws.Range('E27') '  +1.35%  '
Set-TextValue $ws.Range('E28') '  -0.05%  '
Set-TextValue $ws.Range('D29') '7.53'
Set-TextValue $ws.Range('E29') '  +7.12%  '
Set-TextValue $ws.Range('D30') '2.70'
Set-TextValue $ws.Range('E30') '  +2.07%  '
Set-TextValue $ws.Range('B31') 'ImmutableX'
Set-TextValue $ws.Range('C31') 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Range('D31') '2.23'
Set-TextValue $ws.Range('E31') '  +8.30%  '
Set-TextValue $ws.Range('B32') 'FirstDigitalUSD'
Set-TextValue $ws.Range('C32') 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue $ws.Range('D32') '1.00'
Set-TextValue $ws.Range('E32') '  -0.03%  '
Set-TextValue $ws.Range('D33') '27.75'
Set-TextValue $ws.Range('E33') '  +1.81%  '
Set-TextValue $ws.Range('D34') '0.111'
Set-TextValue $ws.Range('E34') '  +3.09%  '
Set-TextValue $ws.Range('D35') '0.0₃0864'
Set-TextValue $ws.Range('E35') '  +6.64%  '
Set-TextValue $ws.Range('D36') '1.05'
Set-TextValue $ws.Range('E36') '  +2.35%  '
Set-TextValue $ws.Range('D37') '5.93'
Set-TextValue $ws.Range('E37') '  +3.15%  '
Set-TextValue $ws.Range('D38') '3.20'
Set-TextValue $ws.Range('E38') '  +11.87%  '
Set-TextValue $ws.Range('B39') 'Stacks'
Set-TextValue $ws.Range('C39') 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws.Range('D39') '2.09'
Set-TextValue $ws.Range('E39') '  +0.24%  '
Set-TextValue $ws.Range('B40') 'OKB'
Set-TextValue $ws.Range('C40') 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws.Range('D40') '50.72'
Set-TextValue $ws.Range('E40') '  +1.00%  '
Set-TextValue $ws.Range('D41') '9.13'
Set-TextValue $ws.Range('E41') '  -0.94%  '
Set-TextValue $ws.Range('D42') '0.126'
Set-TextValue $ws.Range('E42') '  +2.46%  '
Set-TextValue $ws.Range('D43') '0.307'
Set-TextValue $ws.Range('E43') '  +15.84%  '
Set-TextValue $ws.Range('D44') '42.32'
Set-TextValue $ws.Range('E44') '  +9.04%  '
Set-TextValue $ws.Range('D45') '395.59'
Set-TextValue $ws.Range('E45') '  -1.43%  '
Set-TextValue $ws.Range('D46') '0.0360'
Set-TextValue $ws.Range('E46') '  +2.11%  '
Set-TextValue $ws.Range('D47') '2.742.94'
Set-TextValue $ws.Range('E47') '  +0.90%  '
Set-TextValue $ws.Range('D48') '131.71'
Set-TextValue $ws.Range('E48') '  -0.88%  '
Set-TextValue $ws.Range('E49') '  +0.07%  '
Set-TextValue $ws.Range('D50') '2.22'
Set-TextValue $ws.Range('E50') '  +2.62%  '
Set-TextValue $ws.Range('D51') '0.109'
Set-TextValue $ws.Range('E51') '  +1.08%  '
